# Refactor currency conversion sheet: split the single "foreign_amount"
# column into explicit source/target amounts and add a target_fees column.
#
# Old layout (currency_conversions):
#   A: date | B: foreign_amount | C: source_fees | D: source_currency | E: target_currency | F: comment
# New layout:
#   A: date | B: source_amount | C: source_fees | D: source_currency | E: target_amount | F: target_fees | G: target_currency | H: comment

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# Insert two new columns right before the (old) target_currency column (E),
# pushing target_currency -> G and comment -> H, and leaving B..D untouched.
$ws.Range("E1:F1").EntireColumn.Insert()

# Rename the old "foreign_amount" header (still in B1) to "source_amount".
$ws.Range("B1").Value = "source_amount"

# Fill in the headers for the two newly inserted columns.
$ws.Range("E1").Value = "target_amount"
$ws.Range("F1").Value = "target_fees"

# Populate the new target_amount / target_fees values for each data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = -1
    $ws.Cells.Item($r, 6).Value = 0
}

# Make currency_conversions the active sheet/tab (was money_transfers before).
$ws.Activate()
